# Update "想去人数" (F column) values across the four worksheets to reflect
# the newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 353
$ws1.Range("F4").Value = 404
$ws1.Range("F5").Value = 304
$ws1.Range("F8").Value = 815
$ws1.Range("F9").Value = 1578
$ws1.Range("F10").Value = 5991
$ws1.Range("F12").Value = 1711
$ws1.Range("F13").Value = 430
$ws1.Range("F14").Value = 5785
$ws1.Range("F15").Value = 5785
$ws1.Range("F16").Value = 109
$ws1.Range("F18").Value = 147
$ws1.Range("F19").Value = 93
$ws1.Range("F20").Value = 1619
$ws1.Range("F21").Value = 835
$ws1.Range("F23").Value = 134
$ws1.Range("F24").Value = 1301
$ws1.Range("F25").Value = 709
$ws1.Range("F26").Value = 216
$ws1.Range("F28").Value = 2
$ws1.Range("F31").Value = 3847

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 292
$ws2.Range("F5").Value = 147
$ws2.Range("F8").Value = 363
$ws2.Range("F13").Value = 18
$ws2.Range("F20").Value = 59

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 9485
$ws3.Range("F3").Value = 2207
$ws3.Range("F4").Value = 586
$ws3.Range("F5").Value = 142

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 9485
$ws4.Range("F3").Value = 2207
$ws4.Range("F4").Value = 586
$ws4.Range("F5").Value = 353
$ws4.Range("F6").Value = 404
$ws4.Range("F7").Value = 309
$ws4.Range("F11").Value = 292
$ws4.Range("F12").Value = 815
$ws4.Range("F13").Value = 142
$ws4.Range("F14").Value = 1578
$ws4.Range("F15").Value = 5991
$ws4.Range("F17").Value = 363
$ws4.Range("F18").Value = 1711
$ws4.Range("F21").Value = 430
$ws4.Range("F24").Value = 5785
$ws4.Range("F25").Value = 5785
$ws4.Range("F26").Value = 109
$ws4.Range("F28").Value = 147
$ws4.Range("F29").Value = 93
$ws4.Range("F30").Value = 1619
$ws4.Range("F31").Value = 835
$ws4.Range("F33").Value = 134
$ws4.Range("F34").Value = 1301
$ws4.Range("F35").Value = 709
$ws4.Range("F36").Value = 216
$ws4.Range("F38").Value = 18
$ws4.Range("F39").Value = 2
$ws4.Range("F46").Value = 3847
$ws4.Range("F48").Value = 59

$wb.Save()
